# Apply the "7_gunluk_vardiya_plani" shift-plan update:
#   - rewrite the per-hour headcount numbers in rows 4-10 (columns L..Y)
#   - move the active-cell selection to T17 (cosmetic, matches author's last click)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "L4" = 1;  "M4" = 1;  "N4" = 3;  "O4" = 4;  "P4" = 7;
    "S4" = 6;  "T4" = 6;  "U4" = 6;  "V4" = 6;  "W4" = 4;  "X4" = 3;  "Y4" = 3;

    "M5" = 2;  "N5" = 4;  "O5" = 4;  "P5" = 8;  "Q5" = 8;
    "R5" = 8;  "S5" = 8;  "T5" = 6;  "X5" = 4;

    "L6" = 2;  "M6" = 2;  "N6" = 4;  "O6" = 4;  "P6" = 8;
    "Q6" = 8;  "R6" = 8;  "S6" = 8;  "T6" = 6;  "U6" = 6;  "V6" = 4;  "W6" = 4;  "X6" = 4;

    "L7" = 2;  "M7" = 2;  "N7" = 4;  "O7" = 4;  "Q7" = 8;
    "R7" = 8;  "T7" = 6;  "U7" = 6;  "X7" = 4;

    "L8" = 2;  "M8" = 2;  "N8" = 4;  "O8" = 5;  "P8" = 8;
    "Q8" = 8;  "S8" = 8;  "T8" = 6;  "U8" = 6;  "V8" = 6;  "W8" = 4;  "X8" = 3;

    "N9" = 5;  "O9" = 6;  "P9" = 9;  "Q9" = 9;  "R9" = 9;
    "S9" = 8;  "T9" = 8;  "U9" = 6;  "V9" = 6;  "W9" = 4;  "X9" = 3;

    "M10" = 3; "O10" = 6; "P10" = 9; "Q10" = 9; "R10" = 9;
    "S10" = 8; "T10" = 8; "U10" = 6; "V10" = 6; "W10" = 4; "X10" = 3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Restore the selection/active cell left behind by the author
$ws.Range("T17").Select()
